$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.024.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5910"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2751"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06781"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07499"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.808.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.675"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6237"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009422"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "74.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "28.756.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.435"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -8.76%  "

$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.768"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "154.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1271"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.775"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06408"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.403"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.431"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.718"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.91%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.674"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.682"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.049"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.531"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6303"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.746"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.59%  "

$ws.Range("E38").Value = "  -2.72%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01699"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.130.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.25%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8681"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.95%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.971.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.61%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.07%  "

$ws.Range("E46").Value = "  -3.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.572"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05464"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.69%  "

$ws.Range("E49").Value = "  -1.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.249"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.46%  "

$ws.Range("E51").Value = "  +0.19%  "

